{"js": "// Update the date heading and the five rows of division problems in the\n// practice-sheet table. Both the date paragraph and the table cells are\n// addressed positionally (by paragraph/row/column index) rather than by\n// searching for their old text, because several of the new values\n// duplicate other cells' OLD values (e.g. \"95\u00f78=\", \"99\u00f78=\", \"96\u00f76=\",\n// \"13\u00f76=\" appear as both a before- and an after- value at different\n// spots), so a text-based find/replace-all could cascade and clobber\n// cells that were only just written.\n\n// 1) Date heading paragraph (first paragraph in the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-02-19 Monday\", \"Replace\");\n\n// 2) The division-problem table: 5 \"content\" rows, each with 5 cells.\n//    (Rows in between are blank spacer rows and are left untouched.)\nconst table = body.tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// New values, keyed by (content-row-index, column-index). Content rows in\n// the table are at row indices 0, 4, 8, 12, 16.\nconst newValues = [\n  [\"90\u00f72=\", \"33\u00f77=\", \"60\u00f78=\", \"19\u00f78=\", \"90\u00f79=\"],\n  [\"86\u00f76=\", \"62\u00f78=\", \"69\u00f74=\", \"95\u00f78=\", \"58\u00f72=\"],\n  [\"49\u00f75=\", \"99\u00f78=\", \"81\u00f75=\", \"87\u00f78=\", \"41\u00f74=\"],\n  [\"94\u00f75=\", \"97\u00f73=\", \"13\u00f76=\", \"94\u00f76=\", \"49\u00f79=\"],\n  [\"96\u00f76=\", \"80\u00f73=\", \"79\u00f76=\", \"47\u00f74=\", \"66\u00f76=\"],\n];\nconst contentRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < contentRowIndexes.length; i++) {\n  const rowIdx = contentRowIndexes[i];\n  const cells = rows.items[rowIdx].cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < newValues[i].length; c++) {\n    cells.items[c].value = newValues[i][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the five rows of division problems in the\n# practice-sheet table. Everything is addressed positionally (paragraph\n# index / table row+column index) rather than via Find/Replace of the old\n# literal text: several of the NEW values are identical to OTHER cells'\n# OLD values (e.g. \"95\u00f78=\", \"99\u00f78=\", \"96\u00f76=\", \"13\u00f76=\" each show up as both\n# a before- and an after- value at different spots), so a blind\n# Find-and-ReplaceAll could cascade and corrupt a cell that was only just\n# rewritten. Writing straight to each cell's Range.Text avoids that.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-02-19 Monday\"\n\n# 2) The division-problem table: 5 \"content\" rows (1, 5, 9, 13, 17 in\n#    Word's 1-based Cell() indexing), each with 5 cells. The rows in\n#    between are blank spacer rows and are left untouched.\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"90\u00f72=\", \"33\u00f77=\", \"60\u00f78=\", \"19\u00f78=\", \"90\u00f79=\"),\n    @(\"86\u00f76=\", \"62\u00f78=\", \"69\u00f74=\", \"95\u00f78=\", \"58\u00f72=\"),\n    @(\"49\u00f75=\", \"99\u00f78=\", \"81\u00f75=\", \"87\u00f78=\", \"41\u00f74=\"),\n    @(\"94\u00f75=\", \"97\u00f73=\", \"13\u00f76=\", \"94\u00f76=\", \"49\u00f79=\"),\n    @(\"96\u00f76=\", \"80\u00f73=\", \"79\u00f76=\", \"47\u00f74=\", \"66\u00f76=\")\n)\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $contentRows.Length; $i++) {\n    $row = $contentRows[$i]\n    $values = $newValues[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($row, $c).Range.Text = $values[$c - 1]\n    }\n}\n"}
